$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking price strings
# (e.g. "1.008") are not silently converted to numbers by Excel.
# (D38 is excluded since its price value is not being modified;
# applied as two contiguous ranges since multi-area NumberFormat
# assignment only reliably affects the first area.)
$ws.Range("D2:D37").NumberFormat = "@"
$ws.Range("D39:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.360.24'
$ws.Range("E2").Value = '  -1.20%  '
$ws.Range("D3").Value = '1.872.10'
$ws.Range("E3").Value = '  -1.76%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  -2.40%  '
$ws.Range("D5").Value = '314.42'
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("D7").Value = '0.5089'
$ws.Range("E7").Value = '  -2.09%  '
$ws.Range("D8").Value = '0.3937'
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("D9").Value = '0.08374'
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").Value = '1.106'
$ws.Range("E10").Value = '  -2.67%  '
$ws.Range("D11").Value = '6.233'
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("D12").Value = '1.870.94'
$ws.Range("E12").Value = '  -2.59%  '
$ws.Range("D13").Value = '20.42'
$ws.Range("E13").Value = '  -1.41%  '
$ws.Range("D14").Value = '7.239'
$ws.Range("E14").Value = '  -1.16%  '
$ws.Range("D15").Value = '1.007'
$ws.Range("E15").Value = '  -2.78%  '
$ws.Range("D16").Value = '0.00001102'
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("D17").Value = '90.82'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = '0.06694'
$ws.Range("E18").Value = '  -1.83%  '
$ws.Range("D19").Value = '17.64'
$ws.Range("E19").Value = '  -2.03%  '
$ws.Range("D20").Value = '1.007'
$ws.Range("E20").Value = '  -2.26%  '
$ws.Range("D21").Value = '5.934'
$ws.Range("E21").Value = '  -2.69%  '
$ws.Range("D22").Value = '28.389.10'
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("D23").Value = '11.08'
$ws.Range("E23").Value = '  -1.60%  '
$ws.Range("D24").Value = '2.252'
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("D25").Value = '2.089.89'
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").Value = '160.56'
$ws.Range("E26").Value = '  -1.46%  '
$ws.Range("D27").Value = '20.62'
$ws.Range("E27").Value = '  -1.99%  '
$ws.Range("D28").Value = '2.369'
$ws.Range("E28").Value = '  -3.40%  '
$ws.Range("D29").Value = '126.92'
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").Value = '0.1046'
$ws.Range("E30").Value = '  -1.29%  '
$ws.Range("D31").Value = '1.046'
$ws.Range("E31").Value = '  -0.96%  '
$ws.Range("D32").Value = '5.763'
$ws.Range("E32").Value = '  -3.86%  '
$ws.Range("D33").Value = '3.594'
$ws.Range("E33").Value = '  -2.62%  '
$ws.Range("D34").Value = '0.02429'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").Value = '0.06470'
$ws.Range("E35").Value = '  -2.81%  '
$ws.Range("D36").Value = '0.2179'
$ws.Range("E36").Value = '  -1.87%  '
$ws.Range("D37").Value = '8.865'
$ws.Range("E37").Value = '  -6.50%  '
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("D39").Value = '1.190'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("D40").Value = '5.058'
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").Value = '0.6417'
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("D42").Value = '11.13'
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("D43").Value = '1.007'
$ws.Range("E43").Value = '  -2.20%  '
$ws.Range("D44").Value = '0.6037'
$ws.Range("E44").Value = '  -2.10%  '
$ws.Range("D45").Value = '13.05'
$ws.Range("E45").Value = '  -1.80%  '
$ws.Range("D46").Value = '3.689'
$ws.Range("E46").Value = '  -1.85%  '
$ws.Range("D47").Value = '2.004'
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("D48").Value = '121.73'
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.202'
$ws.Range("E49").Value = '  -3.20%  '
$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").Value = '1.187'
$ws.Range("E50").Value = '  -8.67%  '
$ws.Range("D51").Value = '0.06811'
$ws.Range("E51").Value = '  -2.12%  '
